$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 <- old row 25
$ws.Cells.Item(2, 4).Value = 44446
$ws.Cells.Item(2, 8).Value = "Madrigal"
$ws.Cells.Item(2, 9).Value = "Primera"
$ws.Cells.Item(2, 10).Value = 120
$ws.Cells.Item(2, 11).Value = 16000
$ws.Cells.Item(2, 12).Value = 16000
$ws.Cells.Item(2, 13).Value = 16000
$ws.Cells.Item(2, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(2, 16).Value = 400
$ws.Cells.Item(2, 17).Value = 40

# Row 3 <- old row 40
$ws.Cells.Item(3, 4).Value = 45121
$ws.Cells.Item(3, 8).Value = "Argentina(o)"
$ws.Cells.Item(3, 9).Value = "Primera"
$ws.Cells.Item(3, 10).Value = 120
$ws.Cells.Item(3, 11).Value = 16000
$ws.Cells.Item(3, 12).Value = 16000
$ws.Cells.Item(3, 13).Value = 16000
$ws.Cells.Item(3, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(3, 16).Value = 320
$ws.Cells.Item(3, 17).Value = 50

# Row 4 <- old row 26
$ws.Cells.Item(4, 4).Value = 44161
$ws.Cells.Item(4, 8).Value = "Madrigal"
$ws.Cells.Item(4, 9).Value = "Primera"
$ws.Cells.Item(4, 10).Value = 30
$ws.Cells.Item(4, 11).Value = 11000
$ws.Cells.Item(4, 12).Value = 11000
$ws.Cells.Item(4, 13).Value = 11000
$ws.Cells.Item(4, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(4, 16).Value = 275
$ws.Cells.Item(4, 17).Value = 40

# Row 5 <- old row 37
$ws.Cells.Item(5, 4).Value = 44495
$ws.Cells.Item(5, 8).Value = "Madrigal"
$ws.Cells.Item(5, 9).Value = "Primera"
$ws.Cells.Item(5, 10).Value = 130
$ws.Cells.Item(5, 11).Value = 11000
$ws.Cells.Item(5, 12).Value = 11000
$ws.Cells.Item(5, 13).Value = 11000
$ws.Cells.Item(5, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(5, 16).Value = 275
$ws.Cells.Item(5, 17).Value = 40

# Row 6 <- old row 18
$ws.Cells.Item(6, 4).Value = 44876
$ws.Cells.Item(6, 8).Value = "Española"
$ws.Cells.Item(6, 9).Value = "Primera"
$ws.Cells.Item(6, 10).Value = 200
$ws.Cells.Item(6, 11).Value = 12000
$ws.Cells.Item(6, 12).Value = 12000
$ws.Cells.Item(6, 13).Value = 12000
$ws.Cells.Item(6, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(6, 16).Value = 400
$ws.Cells.Item(6, 17).Value = 30

# Row 7 <- old row 30
$ws.Cells.Item(7, 4).Value = 44789
$ws.Cells.Item(7, 8).Value = "Madrigal"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 80
$ws.Cells.Item(7, 11).Value = 16000
$ws.Cells.Item(7, 12).Value = 16000
$ws.Cells.Item(7, 13).Value = 16000
$ws.Cells.Item(7, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(7, 16).Value = 400
$ws.Cells.Item(7, 17).Value = 40

# Row 8 <- old row 11
$ws.Cells.Item(8, 4).Value = 44176
$ws.Cells.Item(8, 8).Value = "Madrigal"
$ws.Cells.Item(8, 9).Value = "Primera"
$ws.Cells.Item(8, 10).Value = 80
$ws.Cells.Item(8, 11).Value = 11000
$ws.Cells.Item(8, 12).Value = 11000
$ws.Cells.Item(8, 13).Value = 11000
$ws.Cells.Item(8, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(8, 16).Value = 275
$ws.Cells.Item(8, 17).Value = 40

# Row 9 <- old row 16
$ws.Cells.Item(9, 4).Value = 44848
$ws.Cells.Item(9, 8).Value = "Española"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 200
$ws.Cells.Item(9, 11).Value = 10000
$ws.Cells.Item(9, 12).Value = 10000
$ws.Cells.Item(9, 13).Value = 10000
$ws.Cells.Item(9, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(9, 16).Value = 333
$ws.Cells.Item(9, 17).Value = 30

# Row 10 <- old row 34
$ws.Cells.Item(10, 4).Value = 44838
$ws.Cells.Item(10, 8).Value = "Española"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 150
$ws.Cells.Item(10, 11).Value = 12000
$ws.Cells.Item(10, 12).Value = 12000
$ws.Cells.Item(10, 13).Value = 12000
$ws.Cells.Item(10, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(10, 16).Value = 400
$ws.Cells.Item(10, 17).Value = 30

# Row 11 <- old row 33
$ws.Cells.Item(11, 4).Value = 44421
$ws.Cells.Item(11, 8).Value = "Española"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 80
$ws.Cells.Item(11, 11).Value = 16500
$ws.Cells.Item(11, 12).Value = 16500
$ws.Cells.Item(11, 13).Value = 16500
$ws.Cells.Item(11, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(11, 16).Value = 550
$ws.Cells.Item(11, 17).Value = 30

# Row 12 <- old row 7
$ws.Cells.Item(12, 4).Value = 44810
$ws.Cells.Item(12, 8).Value = "Madrigal"
$ws.Cells.Item(12, 9).Value = "Primera"
$ws.Cells.Item(12, 10).Value = 100
$ws.Cells.Item(12, 11).Value = 13500
$ws.Cells.Item(12, 12).Value = 14000
$ws.Cells.Item(12, 13).Value = 13750
$ws.Cells.Item(12, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(12, 16).Value = 344
$ws.Cells.Item(12, 17).Value = 40

# Row 13 <- old row 2
$ws.Cells.Item(13, 4).Value = 44799
$ws.Cells.Item(13, 8).Value = "Madrigal"
$ws.Cells.Item(13, 9).Value = "Primera"
$ws.Cells.Item(13, 10).Value = 70
$ws.Cells.Item(13, 11).Value = 15000
$ws.Cells.Item(13, 12).Value = 15000
$ws.Cells.Item(13, 13).Value = 15000
$ws.Cells.Item(13, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(13, 16).Value = 375
$ws.Cells.Item(13, 17).Value = 40

# Row 14 <- old row 15
$ws.Cells.Item(14, 4).Value = 44873
$ws.Cells.Item(14, 8).Value = "Española"
$ws.Cells.Item(14, 9).Value = "Primera"
$ws.Cells.Item(14, 10).Value = 150
$ws.Cells.Item(14, 11).Value = 12000
$ws.Cells.Item(14, 12).Value = 12000
$ws.Cells.Item(14, 13).Value = 12000
$ws.Cells.Item(14, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(14, 16).Value = 400
$ws.Cells.Item(14, 17).Value = 30

# Row 15 <- old row 4
$ws.Cells.Item(15, 4).Value = 44390
$ws.Cells.Item(15, 8).Value = "Española"
$ws.Cells.Item(15, 9).Value = "Primera"
$ws.Cells.Item(15, 10).Value = 80
$ws.Cells.Item(15, 11).Value = 16000
$ws.Cells.Item(15, 12).Value = 16000
$ws.Cells.Item(15, 13).Value = 16000
$ws.Cells.Item(15, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(15, 16).Value = 533
$ws.Cells.Item(15, 17).Value = 30

# Row 16 <- old row 35
$ws.Cells.Item(16, 4).Value = 44757
$ws.Cells.Item(16, 8).Value = "Argentina(o)"
$ws.Cells.Item(16, 9).Value = "Primera"
$ws.Cells.Item(16, 10).Value = 80
$ws.Cells.Item(16, 11).Value = 18000
$ws.Cells.Item(16, 12).Value = 18000
$ws.Cells.Item(16, 13).Value = 18000
$ws.Cells.Item(16, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(16, 16).Value = 360
$ws.Cells.Item(16, 17).Value = 50

# Row 17 <- old row 32
$ws.Cells.Item(17, 4).Value = 44831
$ws.Cells.Item(17, 8).Value = "Madrigal"
$ws.Cells.Item(17, 9).Value = "Primera"
$ws.Cells.Item(17, 10).Value = 180
$ws.Cells.Item(17, 11).Value = 12000
$ws.Cells.Item(17, 12).Value = 13000
$ws.Cells.Item(17, 13).Value = 12444
$ws.Cells.Item(17, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(17, 16).Value = 311
$ws.Cells.Item(17, 17).Value = 40

# Row 18 <- old row 14
$ws.Cells.Item(18, 4).Value = 44778
$ws.Cells.Item(18, 8).Value = "Madrigal"
$ws.Cells.Item(18, 9).Value = "Primera"
$ws.Cells.Item(18, 10).Value = 160
$ws.Cells.Item(18, 11).Value = 15000
$ws.Cells.Item(18, 12).Value = 16000
$ws.Cells.Item(18, 13).Value = 15500
$ws.Cells.Item(18, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(18, 16).Value = 388
$ws.Cells.Item(18, 17).Value = 40

# Row 19 <- old row 29
$ws.Cells.Item(19, 4).Value = 44775
$ws.Cells.Item(19, 8).Value = "Madrigal"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 100
$ws.Cells.Item(19, 11).Value = 16000
$ws.Cells.Item(19, 12).Value = 17000
$ws.Cells.Item(19, 13).Value = 16500
$ws.Cells.Item(19, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(19, 16).Value = 412
$ws.Cells.Item(19, 17).Value = 40

# Row 20 <- old row 24
$ws.Cells.Item(20, 4).Value = 44827
$ws.Cells.Item(20, 8).Value = "Madrigal"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 100
$ws.Cells.Item(20, 11).Value = 12000
$ws.Cells.Item(20, 12).Value = 12000
$ws.Cells.Item(20, 13).Value = 12000
$ws.Cells.Item(20, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(20, 16).Value = 300
$ws.Cells.Item(20, 17).Value = 40

# Row 21 <- old row 28
$ws.Cells.Item(21, 4).Value = 45107
$ws.Cells.Item(21, 8).Value = "Argentina(o)"
$ws.Cells.Item(21, 9).Value = "Primera"
$ws.Cells.Item(21, 10).Value = 60
$ws.Cells.Item(21, 11).Value = 20000
$ws.Cells.Item(21, 12).Value = 20000
$ws.Cells.Item(21, 13).Value = 20000
$ws.Cells.Item(21, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(21, 16).Value = 400
$ws.Cells.Item(21, 17).Value = 50

# Row 22 <- old row 23
$ws.Cells.Item(22, 4).Value = 44488
$ws.Cells.Item(22, 8).Value = "Madrigal"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 120
$ws.Cells.Item(22, 11).Value = 12000
$ws.Cells.Item(22, 12).Value = 12000
$ws.Cells.Item(22, 13).Value = 12000
$ws.Cells.Item(22, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(22, 16).Value = 300
$ws.Cells.Item(22, 17).Value = 40

# Row 23 <- old row 38
$ws.Cells.Item(23, 4).Value = 44859
$ws.Cells.Item(23, 8).Value = "Española"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 150
$ws.Cells.Item(23, 11).Value = 12000
$ws.Cells.Item(23, 12).Value = 12000
$ws.Cells.Item(23, 13).Value = 12000
$ws.Cells.Item(23, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(23, 16).Value = 400
$ws.Cells.Item(23, 17).Value = 30

# Row 24 <- old row 10
$ws.Cells.Item(24, 4).Value = 44400
$ws.Cells.Item(24, 8).Value = "Española"
$ws.Cells.Item(24, 9).Value = "Primera"
$ws.Cells.Item(24, 10).Value = 70
$ws.Cells.Item(24, 11).Value = 15000
$ws.Cells.Item(24, 12).Value = 15000
$ws.Cells.Item(24, 13).Value = 15000
$ws.Cells.Item(24, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(24, 16).Value = 500
$ws.Cells.Item(24, 17).Value = 30

# Row 25 <- old row 31
$ws.Cells.Item(25, 4).Value = 44841
$ws.Cells.Item(25, 8).Value = "Madrigal"
$ws.Cells.Item(25, 9).Value = "Primera"
$ws.Cells.Item(25, 10).Value = 150
$ws.Cells.Item(25, 11).Value = 11000
$ws.Cells.Item(25, 12).Value = 12000
$ws.Cells.Item(25, 13).Value = 11533
$ws.Cells.Item(25, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(25, 16).Value = 384
$ws.Cells.Item(25, 17).Value = 30

# Row 26 <- old row 19
$ws.Cells.Item(26, 4).Value = 44162
$ws.Cells.Item(26, 8).Value = "Madrigal"
$ws.Cells.Item(26, 9).Value = "Primera"
$ws.Cells.Item(26, 10).Value = 50
$ws.Cells.Item(26, 11).Value = 10000
$ws.Cells.Item(26, 12).Value = 10000
$ws.Cells.Item(26, 13).Value = 10000
$ws.Cells.Item(26, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(26, 16).Value = 250
$ws.Cells.Item(26, 17).Value = 40

# Row 27 <- old row 3
$ws.Cells.Item(27, 4).Value = 44855
$ws.Cells.Item(27, 8).Value = "Madrigal"
$ws.Cells.Item(27, 9).Value = "Primera"
$ws.Cells.Item(27, 10).Value = 160
$ws.Cells.Item(27, 11).Value = 10000
$ws.Cells.Item(27, 12).Value = 10000
$ws.Cells.Item(27, 13).Value = 10000
$ws.Cells.Item(27, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(27, 16).Value = 333
$ws.Cells.Item(27, 17).Value = 30

# Row 28 <- old row 27
$ws.Cells.Item(28, 4).Value = 44386
$ws.Cells.Item(28, 8).Value = "Española"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 30
$ws.Cells.Item(28, 11).Value = 15000
$ws.Cells.Item(28, 12).Value = 15000
$ws.Cells.Item(28, 13).Value = 15000
$ws.Cells.Item(28, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(28, 16).Value = 500
$ws.Cells.Item(28, 17).Value = 30

# Row 29 <- old row 20
$ws.Cells.Item(29, 4).Value = 44782
$ws.Cells.Item(29, 8).Value = "Madrigal"
$ws.Cells.Item(29, 9).Value = "Primera"
$ws.Cells.Item(29, 10).Value = 90
$ws.Cells.Item(29, 11).Value = 15000
$ws.Cells.Item(29, 12).Value = 15000
$ws.Cells.Item(29, 13).Value = 15000
$ws.Cells.Item(29, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(29, 16).Value = 375
$ws.Cells.Item(29, 17).Value = 40

# Row 30 <- old row 21
$ws.Cells.Item(30, 4).Value = 44832
$ws.Cells.Item(30, 8).Value = "Madrigal"
$ws.Cells.Item(30, 9).Value = "Primera"
$ws.Cells.Item(30, 10).Value = 120
$ws.Cells.Item(30, 11).Value = 12000
$ws.Cells.Item(30, 12).Value = 13000
$ws.Cells.Item(30, 13).Value = 12500
$ws.Cells.Item(30, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(30, 16).Value = 312
$ws.Cells.Item(30, 17).Value = 40

# Row 31 <- old row 36
$ws.Cells.Item(31, 4).Value = 44771
$ws.Cells.Item(31, 8).Value = "Madrigal"
$ws.Cells.Item(31, 9).Value = "Primera"
$ws.Cells.Item(31, 10).Value = 90
$ws.Cells.Item(31, 11).Value = 16000
$ws.Cells.Item(31, 12).Value = 16000
$ws.Cells.Item(31, 13).Value = 16000
$ws.Cells.Item(31, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(31, 16).Value = 400
$ws.Cells.Item(31, 17).Value = 40

# Row 32 <- old row 12
$ws.Cells.Item(32, 4).Value = 44418
$ws.Cells.Item(32, 8).Value = "Española"
$ws.Cells.Item(32, 9).Value = "Primera"
$ws.Cells.Item(32, 10).Value = 80
$ws.Cells.Item(32, 11).Value = 16000
$ws.Cells.Item(32, 12).Value = 16000
$ws.Cells.Item(32, 13).Value = 16000
$ws.Cells.Item(32, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(32, 16).Value = 533
$ws.Cells.Item(32, 17).Value = 30

# Row 33 <- old row 5
$ws.Cells.Item(33, 4).Value = 44803
$ws.Cells.Item(33, 8).Value = "Madrigal"
$ws.Cells.Item(33, 9).Value = "Primera"
$ws.Cells.Item(33, 10).Value = 100
$ws.Cells.Item(33, 11).Value = 14000
$ws.Cells.Item(33, 12).Value = 15000
$ws.Cells.Item(33, 13).Value = 14500
$ws.Cells.Item(33, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(33, 16).Value = 362
$ws.Cells.Item(33, 17).Value = 40

# Row 34 <- old row 8
$ws.Cells.Item(34, 4).Value = 44484
$ws.Cells.Item(34, 8).Value = "Madrigal"
$ws.Cells.Item(34, 9).Value = "Primera"
$ws.Cells.Item(34, 10).Value = 110
$ws.Cells.Item(34, 11).Value = 11000
$ws.Cells.Item(34, 12).Value = 11000
$ws.Cells.Item(34, 13).Value = 11000
$ws.Cells.Item(34, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(34, 16).Value = 220
$ws.Cells.Item(34, 17).Value = 50

# Row 35 <- old row 9
$ws.Cells.Item(35, 4).Value = 44166
$ws.Cells.Item(35, 8).Value = "Madrigal"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 80
$ws.Cells.Item(35, 11).Value = 10000
$ws.Cells.Item(35, 12).Value = 10000
$ws.Cells.Item(35, 13).Value = 10000
$ws.Cells.Item(35, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(35, 16).Value = 250
$ws.Cells.Item(35, 17).Value = 40

# Row 36 <- old row 17
$ws.Cells.Item(36, 4).Value = 44491
$ws.Cells.Item(36, 8).Value = "Madrigal"
$ws.Cells.Item(36, 9).Value = "Primera"
$ws.Cells.Item(36, 10).Value = 200
$ws.Cells.Item(36, 11).Value = 11000
$ws.Cells.Item(36, 12).Value = 11000
$ws.Cells.Item(36, 13).Value = 11000
$ws.Cells.Item(36, 14).Value = "$/caja 40 unidades"
$ws.Cells.Item(36, 16).Value = 275
$ws.Cells.Item(36, 17).Value = 40

# Row 37 <- old row 13
$ws.Cells.Item(37, 4).Value = 44407
$ws.Cells.Item(37, 8).Value = "Española"
$ws.Cells.Item(37, 9).Value = "Primera"
$ws.Cells.Item(37, 10).Value = 100
$ws.Cells.Item(37, 11).Value = 18000
$ws.Cells.Item(37, 12).Value = 18000
$ws.Cells.Item(37, 13).Value = 18000
$ws.Cells.Item(37, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(37, 16).Value = 600
$ws.Cells.Item(37, 17).Value = 30

# Row 38 <- old row 39
$ws.Cells.Item(38, 4).Value = 44481
$ws.Cells.Item(38, 8).Value = "Madrigal"
$ws.Cells.Item(38, 9).Value = "Segunda"
$ws.Cells.Item(38, 10).Value = 120
$ws.Cells.Item(38, 11).Value = 11000
$ws.Cells.Item(38, 12).Value = 11000
$ws.Cells.Item(38, 13).Value = 11000
$ws.Cells.Item(38, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(38, 16).Value = 220
$ws.Cells.Item(38, 17).Value = 50

# Row 39 <- old row 6
$ws.Cells.Item(39, 4).Value = 44845
$ws.Cells.Item(39, 8).Value = "Española"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 120
$ws.Cells.Item(39, 11).Value = 10000
$ws.Cells.Item(39, 12).Value = 10000
$ws.Cells.Item(39, 13).Value = 10000
$ws.Cells.Item(39, 14).Value = "$/caja 30 unidades"
$ws.Cells.Item(39, 16).Value = 333
$ws.Cells.Item(39, 17).Value = 30

# Row 40 <- old row 22
$ws.Cells.Item(40, 4).Value = 44806
$ws.Cells.Item(40, 8).Value = "Madrigal"
$ws.Cells.Item(40, 9).Value = "Segunda"
$ws.Cells.Item(40, 10).Value = 80
$ws.Cells.Item(40, 11).Value = 13000
$ws.Cells.Item(40, 12).Value = 13000
$ws.Cells.Item(40, 13).Value = 13000
$ws.Cells.Item(40, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(40, 16).Value = 260
$ws.Cells.Item(40, 17).Value = 50
